$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4112
$ws.Range("I74").Value = 3985
$ws.Range("J74").Value = 4158.1816
$ws.Range("K74").Value = 3985
$ws.Range("L74").Value = 4158.1816
$ws.Range("M74").Value = -3049
$ws.Range("N74").Value = -6030.1816
$ws.Range("H77").Value = 4112
$ws.Range("I77").Value = 3985
$ws.Range("J77").Value = 4158.1816
$ws.Range("K77").Value = 19925
$ws.Range("L77").Value = 20790.908
$ws.Range("M77").Value = -15245
$ws.Range("N77").Value = -30150.908
$ws.Range("H112").Value = 2837.7144
$ws.Range("I112").Value = 610
$ws.Range("J112").Value = 3394.6428
$ws.Range("K112").Value = 1830
$ws.Range("L112").Value = 10183.9284
$ws.Range("M112").Value = -722
$ws.Range("N112").Value = -12399.9284
$ws.Range("H116").Value = 3270.5
$ws.Range("I116").Value = 2978.7
$ws.Range("K116").Value = 2978.7
$ws.Range("M116").Value = 463.3000000000002
$ws.Range("H129").Value = 1731.225
$ws.Range("I129").Value = 717.36365
$ws.Range("J129").Value = 2115.7932
$ws.Range("K129").Value = 2152.09095
$ws.Range("L129").Value = 6347.3796
$ws.Range("M129").Value = 2847.90905
$ws.Range("N129").Value = -16347.3796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6023.37
$ws.Range("I32").Value = 4571.3
$ws.Range("J32").Value = 19092
$ws.Range("K32").Value = 4571.3
$ws.Range("L32").Value = 19092
$ws.Range("M32").Value = -4284.3
$ws.Range("N32").Value = -19666
$ws.Range("H57").Value = 26999.5
$ws.Range("I57").Value = 26999.5
$ws.Range("K57").Value = 26999.5
$ws.Range("M57").Value = -26515.5
$ws.Range("H61").Value = 1778.2354
$ws.Range("I61").Value = 1766.1482
$ws.Range("J61").Value = 1791.8334
$ws.Range("K61").Value = 1766.1482
$ws.Range("L61").Value = 1791.8334
$ws.Range("M61").Value = -1554.1482
$ws.Range("N61").Value = -2215.8334
$ws.Range("H88").Value = 1847.6666
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 1847.2
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 1847.2
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -2659.2
$ws.Range("H91").Value = 1847.6666
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 1847.2
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 1847.2
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -4655.2
$ws.Range("H97").Value = 1588.2727
$ws.Range("I97").Value = 1583.3334
$ws.Range("J97").Value = 1594.2
$ws.Range("K97").Value = 1583.3334
$ws.Range("L97").Value = 1594.2
$ws.Range("M97").Value = -1087.3334
$ws.Range("N97").Value = -2586.2
$ws.Range("H132").Value = 1816548.5
$ws.Range("I132").Value = 4478.61
$ws.Range("J132").Value = 4469936.5
$ws.Range("K132").Value = 13435.83
$ws.Range("L132").Value = 13409809.5
$ws.Range("M132").Value = -10905.83
$ws.Range("N132").Value = -13414869.5
$ws.Range("H136").Value = 1778.2354
$ws.Range("I136").Value = 1766.1482
$ws.Range("J136").Value = 1791.8334
$ws.Range("K136").Value = 5298.444600000001
$ws.Range("L136").Value = 5375.5002
$ws.Range("M136").Value = -2748.444600000001
$ws.Range("N136").Value = -10475.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1960
$ws.Range("I86").Value = 1400
$ws.Range("J86").Value = 2333.3333
$ws.Range("K86").Value = 1400
$ws.Range("L86").Value = 2333.3333
$ws.Range("M86").Value = -277
$ws.Range("N86").Value = -4579.3333
$ws.Range("H89").Value = 1960
$ws.Range("I89").Value = 1400
$ws.Range("J89").Value = 2333.3333
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 11666.6665
$ws.Range("M89").Value = -1384
$ws.Range("N89").Value = -22898.6665
$ws.Range("H94").Value = 1665.3334
$ws.Range("I94").Value = 1674.4
$ws.Range("K94").Value = 1674.4
$ws.Range("M94").Value = -1223.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1882.14
$ws.Range("I31").Value = 1123.0518
$ws.Range("J31").Value = 2930.4048
$ws.Range("K31").Value = 1123.0518
$ws.Range("L31").Value = 2930.4048
$ws.Range("M31").Value = -828.0518
$ws.Range("N31").Value = -3520.4048
$ws.Range("H34").Value = 1882.14
$ws.Range("I34").Value = 1123.0518
$ws.Range("J34").Value = 2930.4048
$ws.Range("K34").Value = 1123.0518
$ws.Range("L34").Value = 2930.4048
$ws.Range("M34").Value = -921.0518
$ws.Range("N34").Value = -3334.4048
$ws.Range("H58").Value = 10179.429
$ws.Range("I58").Value = 5668
$ws.Range("J58").Value = 18300
$ws.Range("K58").Value = 5668
$ws.Range("L58").Value = 18300
$ws.Range("M58").Value = -5465
$ws.Range("N58").Value = -18706
$ws.Range("H107").Value = 1932.8334
$ws.Range("I107").Value = 424.75
$ws.Range("J107").Value = 2686.875
$ws.Range("K107").Value = 424.75
$ws.Range("L107").Value = 2686.875
$ws.Range("M107").Value = 1495.25
$ws.Range("N107").Value = -6526.875
$ws.Range("H132").Value = 2109.7368
$ws.Range("I132").Value = 1455.8
$ws.Range("J132").Value = 3150.0908
$ws.Range("K132").Value = 4367.4
$ws.Range("L132").Value = 9450.2724
$ws.Range("M132").Value = -1837.4
$ws.Range("N132").Value = -14510.2724
$ws.Range("H136").Value = 10179.429
$ws.Range("I136").Value = 5668
$ws.Range("J136").Value = 18300
$ws.Range("K136").Value = 17004
$ws.Range("L136").Value = 54900
$ws.Range("M136").Value = -14454
$ws.Range("N136").Value = -60000
$ws.Range("H138").Value = 40214.168
$ws.Range("J138").Value = 42053.637
$ws.Range("L138").Value = 42053.637
$ws.Range("N138").Value = -52333.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1928.5714
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 4500
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -4327
$ws.Range("N16").Value = -6346
$ws.Range("H22").Value = 2567.4211
$ws.Range("I22").Value = 2166.8333
$ws.Range("J22").Value = 2752.3076
$ws.Range("K22").Value = 6500.499899999999
$ws.Range("L22").Value = 8256.9228
$ws.Range("M22").Value = -6331.499899999999
$ws.Range("N22").Value = -8594.9228
$ws.Range("H23").Value = 183.92857
$ws.Range("I23").Value = 93
$ws.Range("J23").Value = 220.3
$ws.Range("K23").Value = 279
$ws.Range("L23").Value = 660.9000000000001
$ws.Range("M23").Value = -44
$ws.Range("N23").Value = -1130.9
$ws.Range("H27").Value = 2567.4211
$ws.Range("I27").Value = 2166.8333
$ws.Range("J27").Value = 2752.3076
$ws.Range("K27").Value = 6500.499899999999
$ws.Range("L27").Value = 8256.9228
$ws.Range("M27").Value = -6398.499899999999
$ws.Range("N27").Value = -8460.9228
$ws.Range("H64").Value = 4755.4287
$ws.Range("I64").Value = 818.5
$ws.Range("J64").Value = 10004.667
$ws.Range("K64").Value = 2455.5
$ws.Range("L64").Value = 30014.001
$ws.Range("M64").Value = -2185.5
$ws.Range("N64").Value = -30554.001
$ws.Range("H67").Value = 4755.4287
$ws.Range("I67").Value = 818.5
$ws.Range("J67").Value = 10004.667
$ws.Range("K67").Value = 2455.5
$ws.Range("L67").Value = 30014.001
$ws.Range("M67").Value = -1519.5
$ws.Range("N67").Value = -31886.001
$ws.Range("H98").Value = 1973.8334
$ws.Range("I98").Value = 650.5
$ws.Range("J98").Value = 2635.5
$ws.Range("K98").Value = 1951.5
$ws.Range("L98").Value = 7906.5
$ws.Range("M98").Value = -453.5
$ws.Range("N98").Value = -10902.5
$ws.Range("H122").Value = 2043.1224
$ws.Range("J122").Value = 2970
$ws.Range("L122").Value = 26730
$ws.Range("N122").Value = -31630

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2017.0869
$ws.Range("I68").Value = 1646.875
$ws.Range("J68").Value = 2863.2856
$ws.Range("K68").Value = 1646.875
$ws.Range("L68").Value = 2863.2856
$ws.Range("M68").Value = -897.875
$ws.Range("N68").Value = -4361.2856
$ws.Range("H71").Value = 2017.0869
$ws.Range("I71").Value = 1646.875
$ws.Range("J71").Value = 2863.2856
$ws.Range("K71").Value = 8234.375
$ws.Range("L71").Value = 14316.428
$ws.Range("M71").Value = -4490.375
$ws.Range("N71").Value = -21804.428
$ws.Range("H93").Value = 3001.4
$ws.Range("I93").Value = 3001.2
$ws.Range("J93").Value = 3001.6
$ws.Range("K93").Value = 3001.2
$ws.Range("L93").Value = 3001.6
$ws.Range("M93").Value = -1753.2
$ws.Range("N93").Value = -5497.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2624611.5
$ws.Range("J5").Value = 2624611.5
$ws.Range("L5").Value = 2624611.5
$ws.Range("N5").Value = -2624835.5
$ws.Range("H132").Value = 1497.8136
$ws.Range("I132").Value = 1236.3695
$ws.Range("J132").Value = 2422.923
$ws.Range("K132").Value = 3709.1085
$ws.Range("L132").Value = 7268.768999999999
$ws.Range("M132").Value = -1179.1085
$ws.Range("N132").Value = -12328.769
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
